$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15
$ws.Cells.Item($row, 1).Value = 39
$ws.Cells.Item($row, 2).Value = "trying squashes"
$ws.Cells.Item($row, 3).Value = "riya-morankar"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "edit1 to main"

# Force the date column to be stored as text (matching the other rows'
# plain-text date cells) instead of being auto-converted to a date serial.
$ws.Cells.Item($row, 6).NumberFormat = "@"
$ws.Cells.Item($row, 6).Value = "2025-06-18"
$ws.Cells.Item($row, 6).Style = "Normal"
